$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("part2")

# Add new row of diary data at the end (row 7)
$ws.Range("A7").Value = 211129
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = "Test 2.9 and do 2.10"

# Update selection to match the target state
$ws.Activate()
$ws.Range("I11").Select()
